# Apply "Dummy data attractions done" edits to the Attractions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attractions")

# --- Update the workbook-level absolute path metadata (mc:AlternateContent / x15ac:absPath) ---
$wb.AbsolutePath = "F:\0001 MY COURSES\0000 Computer Science\COMP 346 - Internet Computing\346-term-project\Github\danglingpointers\data\"

# --- Stash a copy of the existing "text" number format (style index already used by column D) ---
# on a scratch cell far outside the used range, so we can re-apply the SAME style (no new
# style table entries) to each D-column cell after writing a genuine numeric value into it.
$scratch = $ws.Cells.Item(500, 60)
$ws.Cells.Item(2, 4).Copy()
$scratch.PasteSpecial(-4122)

# Row 2: Chicago / Millenium Park
$ws.Cells.Item(2, 2).Value = 'Parks'
$ws.Cells.Item(2, 3).Value = 'Millenium Park'
$dCell = $ws.Cells.Item(2, 4)
$dCell.ClearFormats()
$dCell.Value = 41.882552
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Value = -87.622551

# Row 3: Chicago / Lincoln Park Zoo
$ws.Cells.Item(3, 2).Value = 'Zoo'
$ws.Cells.Item(3, 3).Value = 'Lincoln Park Zoo'
$dCell = $ws.Cells.Item(3, 4)
$dCell.ClearFormats()
$dCell.Value = 41.92089
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(3, 5).Value = -87.632917

# Row 4: Chicago / Navy Pier
$ws.Cells.Item(4, 2).Value = 'Parks'
$ws.Cells.Item(4, 3).Value = 'Navy Pier'
$dCell = $ws.Cells.Item(4, 4)
$dCell.ClearFormats()
$dCell.Value = 41.891642
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(4, 5).Value = -87.605144

# Row 5: Chicago / The Cloud Gate aka the 'Bean'
$ws.Cells.Item(5, 2).Value = 'Landmarks'
$ws.Cells.Item(5, 3).Value = 'The Cloud Gate aka the ''Bean'''
$dCell = $ws.Cells.Item(5, 4)
$dCell.ClearFormats()
$dCell.Value = 41.882657
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(5, 5).Value = -87.623304

# Row 6: Chicago / Fear City Haunted House
$ws.Cells.Item(6, 2).Value = 'Theater'
$ws.Cells.Item(6, 3).Value = 'Fear City Haunted House'
$dCell = $ws.Cells.Item(6, 4)
$dCell.ClearFormats()
$dCell.Value = 42.031595
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(6, 5).Value = -87.779088

# Row 7: Chicago / Museum of Science and Industry
$ws.Cells.Item(7, 2).Value = 'Museum'
$ws.Cells.Item(7, 3).Value = 'Museum of Science and Industry'
$dCell = $ws.Cells.Item(7, 4)
$dCell.ClearFormats()
$dCell.Value = 41.790573
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(7, 5).Value = -87.583066

# Row 8: Chicago / The Field Museum
$ws.Cells.Item(8, 2).Value = 'Museum'
$ws.Cells.Item(8, 3).Value = 'The Field Museum'
$dCell = $ws.Cells.Item(8, 4)
$dCell.ClearFormats()
$dCell.Value = 41.866261
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(8, 5).Value = -87.61698

# Row 9: Chicago / Chicago Cultural Center
$ws.Cells.Item(9, 2).Value = 'Museum'
$ws.Cells.Item(9, 3).Value = 'Chicago Cultural Center'
$dCell = $ws.Cells.Item(9, 4)
$dCell.ClearFormats()
$dCell.Value = 41.883754
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(9, 5).Value = -87.624952

# Row 10: Chicago / Buckingham Fountain
$ws.Cells.Item(10, 2).Value = 'Landmarks'
$ws.Cells.Item(10, 3).Value = 'Buckingham Fountain'
$dCell = $ws.Cells.Item(10, 4)
$dCell.ClearFormats()
$dCell.Value = 41.875794
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(10, 5).Value = -87.618948

# Row 11: Chicago / The Art Institute of Chicago
$ws.Cells.Item(11, 2).Value = 'Museum'
$ws.Cells.Item(11, 3).Value = 'The Art Institute of Chicago'
$dCell = $ws.Cells.Item(11, 4)
$dCell.ClearFormats()
$dCell.Value = 41.879584
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(11, 5).Value = -87.623713

# Row 12: Chicago / Willis Tower
$ws.Cells.Item(12, 2).Value = 'Landmarks'
$ws.Cells.Item(12, 3).Value = 'Willis Tower'
$dCell = $ws.Cells.Item(12, 4)
$dCell.ClearFormats()
$dCell.Value = 41.878876
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(12, 5).Value = -87.635915

# Row 13: Chicago / 360 Chicago
$ws.Cells.Item(13, 2).Value = 'Landmarks'
$ws.Cells.Item(13, 3).Value = '360 Chicago'
$dCell = $ws.Cells.Item(13, 4)
$dCell.ClearFormats()
$dCell.Value = 41.89851
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(13, 5).Value = -87.622765

# Row 14: Chicago / Chicago Architechture Foundation
$ws.Cells.Item(14, 2).Value = 'Museum'
$ws.Cells.Item(14, 3).Value = 'Chicago Architechture Foundation'
$dCell = $ws.Cells.Item(14, 4)
$dCell.ClearFormats()
$dCell.Value = 41.878521
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(14, 5).Value = -87.624853

# Row 15: Chicago / Haunted Sanitarium
$ws.Cells.Item(15, 2).Value = 'Theater'
$ws.Cells.Item(15, 3).Value = 'Haunted Sanitarium'
$dCell = $ws.Cells.Item(15, 4)
$dCell.ClearFormats()
$dCell.Value = 41.927099
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(15, 5).Value = -87.630703

# Row 16: Chicago / Chicago Bulls
$ws.Cells.Item(16, 2).Value = 'Sports'
$ws.Cells.Item(16, 3).Value = 'Chicago Bulls'
$dCell = $ws.Cells.Item(16, 4)
$dCell.ClearFormats()
$dCell.Value = 41.8808159
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(16, 5).Value = -87.7189628

# Row 17: San Francisco / Golden Gate Park
$ws.Cells.Item(17, 2).Value = 'Parks'
$ws.Cells.Item(17, 3).Value = 'Golden Gate Park'
$dCell = $ws.Cells.Item(17, 4)
$dCell.ClearFormats()
$dCell.Value = 37.769421
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(17, 5).Value = -122.486214

# Row 18: San Francisco / San Francisco Zoo
$ws.Cells.Item(18, 2).Value = 'Zoo'
$ws.Cells.Item(18, 3).Value = 'San Francisco Zoo'
$dCell = $ws.Cells.Item(18, 4)
$dCell.ClearFormats()
$dCell.Value = 37.73284
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(18, 5).Value = -122.503065

# Row 19: San Francisco / Conservatory of Flowers
$ws.Cells.Item(19, 2).Value = 'Parks'
$ws.Cells.Item(19, 3).Value = 'Conservatory of Flowers'
$dCell = $ws.Cells.Item(19, 4)
$dCell.ClearFormats()
$dCell.Value = 37.772668
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(19, 5).Value = -122.458758

# Row 20: San Francisco / Palace of Fine Arts
$ws.Cells.Item(20, 2).Value = 'Parks'
$ws.Cells.Item(20, 3).Value = 'Palace of Fine Arts'
$dCell = $ws.Cells.Item(20, 4)
$dCell.ClearFormats()
$dCell.Value = 37.801456
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(20, 5).Value = -122.448053

# Row 21: San Francisco / The San Francisco Dungeon
$ws.Cells.Item(21, 2).Value = 'Theater'
$ws.Cells.Item(21, 3).Value = 'The San Francisco Dungeon'
$dCell = $ws.Cells.Item(21, 4)
$dCell.ClearFormats()
$dCell.Value = 37.808237
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(21, 5).Value = -122.41574

# Row 22: San Francisco / Musee Mecanique
$ws.Cells.Item(22, 2).Value = 'Museum'
$ws.Cells.Item(22, 3).Value = 'Musee Mecanique'
$dCell = $ws.Cells.Item(22, 4)
$dCell.ClearFormats()
$dCell.Value = 37.809305
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(22, 5).Value = -122.416021

# Row 23: San Francisco / California Academy of Sciences
$ws.Cells.Item(23, 2).Value = 'Museum'
$ws.Cells.Item(23, 3).Value = 'California Academy of Sciences'
$dCell = $ws.Cells.Item(23, 4)
$dCell.ClearFormats()
$dCell.Value = 37.769865
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(23, 5).Value = -122.466095

# Row 24: San Francisco / The Walt Disney Family Museum
$ws.Cells.Item(24, 2).Value = 'Museum'
$ws.Cells.Item(24, 3).Value = 'The Walt Disney Family Museum'
$dCell = $ws.Cells.Item(24, 4)
$dCell.ClearFormats()
$dCell.Value = 37.801395
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(24, 5).Value = -122.458661

# Row 25: San Francisco / Golden Gate Bridge
$ws.Cells.Item(25, 2).Value = 'Landmarks'
$ws.Cells.Item(25, 3).Value = 'Golden Gate Bridge'
$dCell = $ws.Cells.Item(25, 4)
$dCell.ClearFormats()
$dCell.Value = 37.819929
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(25, 5).Value = -122.478255

# Row 26: San Francisco / Randall Museum
$ws.Cells.Item(26, 2).Value = 'Museum'
$ws.Cells.Item(26, 3).Value = 'Randall Museum'
$dCell = $ws.Cells.Item(26, 4)
$dCell.ClearFormats()
$dCell.Value = 37.764324
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(26, 5).Value = -122.438408

# Row 27: San Francisco / Alcatraz Island
$ws.Cells.Item(27, 2).Value = 'Landmarks'
$ws.Cells.Item(27, 3).Value = 'Alcatraz Island'
$dCell = $ws.Cells.Item(27, 4)
$dCell.ClearFormats()
$dCell.Value = 37.826978
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(27, 5).Value = -122.422956

# Row 28: San Francisco / Lombard Street
$ws.Cells.Item(28, 2).Value = 'Landmarks'
$ws.Cells.Item(28, 3).Value = 'Lombard Street'
$dCell = $ws.Cells.Item(28, 4)
$dCell.ClearFormats()
$dCell.Value = 37.802139
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(28, 5).Value = -122.41874

# Row 29: San Francisco / Exploratorium
$ws.Cells.Item(29, 2).Value = 'Museum'
$ws.Cells.Item(29, 3).Value = 'Exploratorium'
$dCell = $ws.Cells.Item(29, 4)
$dCell.ClearFormats()
$dCell.Value = 37.800856
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(29, 5).Value = -122.398635

# Row 30: San Francisco / Castro Theatre
$ws.Cells.Item(30, 2).Value = 'Theater'
$ws.Cells.Item(30, 3).Value = 'Castro Theatre'
$dCell = $ws.Cells.Item(30, 4)
$dCell.ClearFormats()
$dCell.Value = 37.762014
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(30, 5).Value = -122.434924

# Row 31: San Francisco / San Francisco Giants
$ws.Cells.Item(31, 2).Value = 'Sports'
$ws.Cells.Item(31, 3).Value = 'San Francisco Giants'
$dCell = $ws.Cells.Item(31, 4)
$dCell.ClearFormats()
$dCell.Value = 37.762014
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(31, 5).Value = -122.434924

# Row 32: Houston / Hermann Park Conservancy
$ws.Cells.Item(32, 2).Value = 'Parks'
$ws.Cells.Item(32, 3).Value = 'Hermann Park Conservancy'
$dCell = $ws.Cells.Item(32, 4)
$dCell.ClearFormats()
$dCell.Value = 29.714875
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(32, 5).Value = -95.389183

# Row 33: Houston / The Houston Zoo
$ws.Cells.Item(33, 2).Value = 'Zoo'
$ws.Cells.Item(33, 3).Value = 'The Houston Zoo'
$dCell = $ws.Cells.Item(33, 4)
$dCell.ClearFormats()
$dCell.Value = 29.71191
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(33, 5).Value = -95.392828

# Row 34: Houston / Menil Park
$ws.Cells.Item(34, 2).Value = 'Parks'
$ws.Cells.Item(34, 3).Value = 'Menil Park'
$dCell = $ws.Cells.Item(34, 4)
$dCell.ClearFormats()
$dCell.Value = 29.737218
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(34, 5).Value = -95.397064

# Row 35: Houston / Discovery Green
$ws.Cells.Item(35, 2).Value = 'Parks'
$ws.Cells.Item(35, 3).Value = 'Discovery Green'
$dCell = $ws.Cells.Item(35, 4)
$dCell.ClearFormats()
$dCell.Value = 29.754132
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(35, 5).Value = -95.35991

# Row 36: Houston / Phobia Haunted Houses
$ws.Cells.Item(36, 2).Value = 'Theater'
$ws.Cells.Item(36, 3).Value = 'Phobia Haunted Houses'
$dCell = $ws.Cells.Item(36, 4)
$dCell.ClearFormats()
$dCell.Value = 29.897554
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(36, 5).Value = -95.595484

# Row 37: Houston / Space Center Houston
$ws.Cells.Item(37, 2).Value = 'Museum'
$ws.Cells.Item(37, 3).Value = 'Space Center Houston'
$dCell = $ws.Cells.Item(37, 4)
$dCell.ClearFormats()
$dCell.Value = 29.550201
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(37, 5).Value = -95.097012

# Row 38: Houston / Art Car Museum
$ws.Cells.Item(38, 2).Value = 'Museum'
$ws.Cells.Item(38, 3).Value = 'Art Car Museum'
$dCell = $ws.Cells.Item(38, 4)
$dCell.ClearFormats()
$dCell.Value = 29.772033
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(38, 5).Value = -95.396867

# Row 39: Houston / The Orange Show
$ws.Cells.Item(39, 2).Value = 'Museum'
$ws.Cells.Item(39, 3).Value = 'The Orange Show'
$dCell = $ws.Cells.Item(39, 4)
$dCell.ClearFormats()
$dCell.Value = 29.717559
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(39, 5).Value = -95.324261

# Row 40: Houston / Buffalo Bayou Park
$ws.Cells.Item(40, 2).Value = 'Parks'
$ws.Cells.Item(40, 3).Value = 'Buffalo Bayou Park'
$dCell = $ws.Cells.Item(40, 4)
$dCell.ClearFormats()
$dCell.Value = 29.761621
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(40, 5).Value = -95.393729

# Row 41: Houston / Houston Museum of Natural Science
$ws.Cells.Item(41, 2).Value = 'Museum'
$ws.Cells.Item(41, 3).Value = 'Houston Museum of Natural Science'
$dCell = $ws.Cells.Item(41, 4)
$dCell.ClearFormats()
$dCell.Value = 29.721819
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(41, 5).Value = -95.389702

# Row 42: Houston / Gerald D. Hines Waterfall
$ws.Cells.Item(42, 2).Value = 'Landmarks'
$ws.Cells.Item(42, 3).Value = 'Gerald D. Hines Waterfall'
$dCell = $ws.Cells.Item(42, 4)
$dCell.ClearFormats()
$dCell.Value = 29.737152
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(42, 5).Value = -95.461068

# Row 43: Houston / Eleanor Tinsley Park
$ws.Cells.Item(43, 2).Value = 'Parks'
$ws.Cells.Item(43, 3).Value = 'Eleanor Tinsley Park'
$dCell = $ws.Cells.Item(43, 4)
$dCell.ClearFormats()
$dCell.Value = 29.761644
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(43, 5).Value = -95.377597

# Row 44: Houston / The Menil Collection
$ws.Cells.Item(44, 2).Value = 'Museum'
$ws.Cells.Item(44, 3).Value = 'The Menil Collection'
$dCell = $ws.Cells.Item(44, 4)
$dCell.ClearFormats()
$dCell.Value = 29.73734
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(44, 5).Value = -95.39851

# Row 45: Houston / Da Camera
$ws.Cells.Item(45, 2).Value = 'Theater'
$ws.Cells.Item(45, 3).Value = 'Da Camera'
$dCell = $ws.Cells.Item(45, 4)
$dCell.ClearFormats()
$dCell.Value = 29.736794
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(45, 5).Value = -95.397309

# Row 46: Houston / Houston Rockets
$ws.Cells.Item(46, 2).Value = 'Sports'
$ws.Cells.Item(46, 3).Value = 'Houston Rockets'
$dCell = $ws.Cells.Item(46, 4)
$dCell.ClearFormats()
$dCell.Value = 29.750767
$scratch.Copy()
$dCell.PasteSpecial(-4122)
$ws.Cells.Item(46, 5).Value = -95.362036

# Clean up the scratch cell so it does not widen the sheet dimension / leave stray data.
$scratch.Clear()

# --- Update the active selection to match the final state (E2) ---
$ws.Range("E2").Select()

Write-Host "edits applied"